# Redo the retry scopes for extraction - new sample record (ezderm) overwrites
# the previous row 2 data. Several columns are no longer populated (CollectionDate,
# Hospital, City, Num of Specimens) while others are newly populated
# (ZipCode, Patient Name, Type Of Procedure 1), and the remaining columns get
# refreshed values for the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that are no longer populated for this record.
$ws.Range("B2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("S2").ClearContents()

# Set the new values, in the same order they first appear left-to-right so the
# shared string table is built up identically to the target workbook.
$ws.Range("A2").Value = "REQUISITION:P2418L0007"
$ws.Range("C2").Value = "Erez A Minka"
$ws.Range("F2").Value = "'76230"
$ws.Range("G2").Value = "940-687-3376"
$ws.Range("H2").Value = "4327 Barnett Road Wichita Falls. TX 763102303"
$ws.Range("I2").Value = "Hutson, Larry"
$ws.Range("J2").Value = "(M/67)"
$ws.Range("K2").Value = "'2024-03-24"
$ws.Range("L2").Value = "Dx 150 Zipper Street Bowie, TX"
$ws.Range("M2").Value = "LAHU0001"
$ws.Range("N2").Value = "LAHU0001"
$ws.Range("O2").Value = "Medicare = Texas"
$ws.Range("P2").Value = "(9A71WX2EA16)"
$ws.Range("Q2").Value = "State Farm Hutson,"
$ws.Range("R2").Value = "(HK4420994343)"
$ws.Range("T2").Value = "Right Ear = Superior Helix 0.5"
$ws.Range("U2").Value = "Neoplasm of uncertain behavior of skin"
$ws.Range("V2").Value = "Biopsy (Tangential (Shave))"
$ws.Range("W2").Value = "Pink papule (Right Ear = Superior Helix)"
$ws.Range("X2").Value = "Right Lateral Forearm ="
$ws.Range("Y2").Value = "Neoplasm of uncertain behavior of skin"
$ws.Range("Z2").Value = "Biopsy (Tangential (Shave))"
$ws.Range("AA2").Value = "Central 0.720.6"
$ws.Range("AB2").Value = "GA"
$ws.Range("AC2").Value = "'303740800"
